$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $text)
    $c = $ws.Range($cellRef)
    $escaped = $text -replace '"', '""'
    $c.Formula = '="' + $escaped + '"'
    $c.Copy()
    $c.PasteSpecial(-4163)
}

Set-TextValue "D2" "61.978.62"
Set-TextValue "E2" "  -0.31%  "
Set-TextValue "D3" "2.459.41"
Set-TextValue "E3" "  -1.50%  "
Set-TextValue "E4" "  +0.10%  "
Set-TextValue "D5" "549.77"
Set-TextValue "E5" "  -0.96%  "
Set-TextValue "D6" "146.53"
Set-TextValue "E6" "  -0.73%  "
Set-TextValue "E7" "  +0.08%  "
Set-TextValue "E8" "  -3.14%  "
Set-TextValue "D9" "2.459.68"
Set-TextValue "E9" "  -1.47%  "
Set-TextValue "E10" "  -1.54%  "
Set-TextValue "E11" "  +0.46%  "
Set-TextValue "D12" "5.43"
Set-TextValue "E12" "  -0.10%  "
Set-TextValue "E13" "  -2.87%  "
Set-TextValue "D14" "26.06"
Set-TextValue "E14" "  -0.96%  "
Set-TextValue "D15" "2.901.96"
Set-TextValue "E15" "  -1.31%  "
Set-TextValue "D16" "0.0000170"
Set-TextValue "E16" "  +2.34%  "
Set-TextValue "D17" "61.842.97"
Set-TextValue "E17" "  -0.27%  "
Set-TextValue "D18" "2.455.77"
Set-TextValue "E18" "  -1.64%  "
Set-TextValue "E19" "  -3.41%  "
Set-TextValue "E20" "  -0.64%  "
Set-TextValue "E21" "  -2.47%  "
Set-TextValue "D22" "320.07"
Set-TextValue "E22" "  -1.38%  "
Set-TextValue "D23" "1.00"
Set-TextValue "E23" "  +0.17%  "
Set-TextValue "E24" "  +7.30%  "
Set-TextValue "D25" "64.08"
Set-TextValue "E25" "  -1.07%  "
Set-TextValue "D26" "0.0₃0982"
Set-TextValue "E26" "  -5.14%  "
Set-TextValue "D27" "2.580.91"
Set-TextValue "E27" "  -2.22%  "
Set-TextValue "E28" "  -0.08%  "
Set-TextValue "E29" "  -2.11%  "
Set-TextValue "D30" "7.86"
Set-TextValue "E30" "  +1.48%  "
Set-TextValue "D31" "531.00"
Set-TextValue "E31" "  -2.27%  "
Set-TextValue "D32" "8.24"
Set-TextValue "E32" "  -3.29%  "
Set-TextValue "E33" "  -3.72%  "
Set-TextValue "E34" "  -1.55%  "
Set-TextValue "D35" "1.63"
Set-TextValue "E35" "  +1.77%  "
Set-TextValue "D36" "5.70"
Set-TextValue "E36" "  -4.07%  "
Set-TextValue "E37" "  +0.30%  "
Set-TextValue "D38" "4.78"
Set-TextValue "E38" "  -2.18%  "
Set-TextValue "E39" "  +0.47%  "
Set-TextValue "D40" "18.24"
Set-TextValue "E40" "  -2.37%  "
Set-TextValue "D41" "1.76"
Set-TextValue "E41" "  +2.38%  "
Set-TextValue "D42" "140.14"
Set-TextValue "E42" "  -4.56%  "
Set-TextValue "E43" "  +0.17%  "
Set-TextValue "D44" "40.40"
Set-TextValue "E45" "  -2.87%  "
Set-TextValue "D46" "144.18"
Set-TextValue "E46" "  -3.45%  "
Set-TextValue "D47" "3.61"
Set-TextValue "E47" "  -0.96%  "
Set-TextValue "D48" "21.37"
Set-TextValue "E48" "  -0.77%  "
Set-TextValue "D49" "0.0529"
Set-TextValue "E49" "  -2.76%  "
Set-TextValue "E50" "  -0.48%  "
Set-TextValue "D51" "0.0934"
Set-TextValue "E51" "  -2.41%  "
